$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 2833.3333
$ws.Range("I51").Value = 2001
$ws.Range("J51").Value = 2937.375
$ws.Range("K51").Value = 2001
$ws.Range("L51").Value = 2937.375
$ws.Range("M51").Value = -1517
$ws.Range("N51").Value = -3905.375
# Row 106
$ws.Range("H106").Value = 3270689.5
$ws.Range("I106").Value = 3586865.8
$ws.Range("K106").Value = 3586865.8
$ws.Range("M106").Value = -3586234.8
# Row 107
$ws.Range("H107").Value = 483459.25
$ws.Range("I107").Value = 529452.0600000001
$ws.Range("K107").Value = 529452.0600000001
$ws.Range("M107").Value = -527532.0600000001
# Row 113
$ws.Range("H113").Value = 86403.25
$ws.Range("I113").Value = 128105
$ws.Range("J113").Value = 2999.75
$ws.Range("K113").Value = 128105
$ws.Range("L113").Value = 2999.75
$ws.Range("M113").Value = -124851
$ws.Range("N113").Value = -9507.75
# Row 129
$ws.Range("H129").Value = 1855.7858
$ws.Range("J129").Value = 1952.6154
$ws.Range("L129").Value = 5857.8462
$ws.Range("N129").Value = -15857.8462
# Row 137
$ws.Range("H137").Value = 18869018
$ws.Range("I137").Value = 25000828
$ws.Range("J137").Value = 1906
$ws.Range("K137").Value = 75002484
$ws.Range("L137").Value = 5718
$ws.Range("M137").Value = -74999934
$ws.Range("N137").Value = -10818

$ws = $wb.Worksheets.Item("ARM")
# Row 9
$ws.Range("H9").Value = 20000
$ws.Range("I9").Value = 50000
$ws.Range("J9").Value = 5000
$ws.Range("K9").Value = 50000
$ws.Range("L9").Value = 5000
$ws.Range("M9").Value = -49830
$ws.Range("N9").Value = -5340
# Row 20
$ws.Range("H20").Value = 20000
$ws.Range("I20").Value = 50000
$ws.Range("J20").Value = 5000
$ws.Range("K20").Value = 50000
$ws.Range("L20").Value = 5000
$ws.Range("M20").Value = -49730
$ws.Range("N20").Value = -5540
# Row 21
$ws.Range("H21").Value = 30406
$ws.Range("I21").Value = 35507.5
$ws.Range("K21").Value = 35507.5
$ws.Range("M21").Value = -35133.5
# Row 22
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
# Row 32
$ws.Range("H32").Value = 24040.285
$ws.Range("I32").Value = 3368.6956
$ws.Range("K32").Value = 3368.6956
$ws.Range("M32").Value = -3081.6956
# Row 74
$ws.Range("H74").Value = 4306.625
$ws.Range("I74").Value = 1019.55884
$ws.Range("J74").Value = 22933.334
$ws.Range("K74").Value = 1019.55884
$ws.Range("L74").Value = 22933.334
$ws.Range("M74").Value = -145.55884
$ws.Range("N74").Value = -24681.334
# Row 77
$ws.Range("H77").Value = 4306.625
$ws.Range("I77").Value = 1019.55884
$ws.Range("J77").Value = 22933.334
$ws.Range("K77").Value = 5097.7942
$ws.Range("L77").Value = 114666.67
$ws.Range("M77").Value = -729.7942000000003
$ws.Range("N77").Value = -123402.67
# Row 102
$ws.Range("H102").Value = 1225
$ws.Range("I102").Value = 1225
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1225
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 397
$ws.Range("N102").ClearContents()
# Row 110
$ws.Range("H110").Value = 654.2
$ws.Range("I110").Value = 561.2
$ws.Range("J110").Value = 933.2
$ws.Range("K110").Value = 561.2
$ws.Range("L110").Value = 933.2
$ws.Range("M110").Value = 1483.8
$ws.Range("N110").Value = -5023.2

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1089.1818
$ws.Range("I107").Value = 986.7778
$ws.Range("J107").Value = 1550
$ws.Range("K107").Value = 986.7778
$ws.Range("L107").Value = 1550
$ws.Range("M107").Value = 933.2222
$ws.Range("N107").Value = -5390
# Row 134
$ws.Range("H134").Value = 2107.1128
$ws.Range("I134").Value = 1250.7819
$ws.Range("J134").Value = 5050.75
$ws.Range("K134").Value = 3752.3457
$ws.Range("L134").Value = 15152.25
$ws.Range("M134").Value = -1217.3457
$ws.Range("N134").Value = -20222.25

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 72742.71000000001
$ws.Range("I16").Value = 125651.5
$ws.Range("K16").Value = 125651.5
$ws.Range("M16").Value = -125364.5
# Row 105
$ws.Range("H105").Value = 903
$ws.Range("I105").Value = 838.5294
$ws.Range("K105").Value = 838.5294
$ws.Range("M105").Value = 908.4706
# Row 107
$ws.Range("H107").Value = 350.14285
$ws.Range("I107").Value = 120.4
$ws.Range("J107").Value = 477.77777
$ws.Range("K107").Value = 120.4
$ws.Range("L107").Value = 477.77777
$ws.Range("M107").Value = 1799.6
$ws.Range("N107").Value = -4317.77777
# Row 113
$ws.Range("H113").Value = 72742.71000000001
$ws.Range("I113").Value = 125651.5
$ws.Range("K113").Value = 125651.5
$ws.Range("M113").Value = -123481.5

$ws = $wb.Worksheets.Item("CUL")
# Row 15
$ws.Range("H15").Value = 366.66666
$ws.Range("I15").Value = 420
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 1260
$ws.Range("L15").Value = 300
$ws.Range("M15").Value = -1120
$ws.Range("N15").Value = -580
# Row 16
$ws.Range("H16").Value = 800.1667
$ws.Range("I16").Value = 700.3333
$ws.Range("K16").Value = 2100.9999
$ws.Range("M16").Value = -1927.9999
# Row 20
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
# Row 22
$ws.Range("H22").Value = 2500
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2500
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 7500
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -7838
# Row 26
$ws.Range("H26").Value = 2950.25
$ws.Range("I26").Value = 400.5
$ws.Range("J26").Value = 5500
$ws.Range("K26").Value = 1201.5
$ws.Range("L26").Value = 16500
$ws.Range("M26").Value = -913.5
$ws.Range("N26").Value = -17076
# Row 27
$ws.Range("H27").Value = 2500
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 2500
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 7500
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -7704
# Row 32
$ws.Range("H32").Value = 10000
$ws.Range("J32").Value = 10000
$ws.Range("L32").Value = 30000
$ws.Range("N32").Value = -30566
# Row 98
$ws.Range("H98").Value = 530
$ws.Range("I98").Value = 420
$ws.Range("K98").Value = 1260
$ws.Range("M98").Value = 238
# Row 113
$ws.Range("H113").Value = 835.95654
$ws.Range("I113").Value = 670
$ws.Range("J113").Value = 924.4666999999999
$ws.Range("K113").Value = 2010
$ws.Range("L113").Value = 2773.4001
$ws.Range("M113").Value = 160
$ws.Range("N113").Value = -7113.4001

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
# Row 113
$ws.Range("H113").Value = 2140.9473
$ws.Range("I113").Value = 1859
$ws.Range("J113").Value = 2346
$ws.Range("K113").Value = 1859
$ws.Range("L113").Value = 2346
$ws.Range("M113").Value = 311
$ws.Range("N113").Value = -6686

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
# Row 22
$ws.Range("H22").Value = 18231.666
$ws.Range("I22").Value = 799.5
$ws.Range("J22").Value = 26947.75
$ws.Range("K22").Value = 799.5
$ws.Range("L22").Value = 26947.75
$ws.Range("M22").Value = -504.5
$ws.Range("N22").Value = -27537.75
# Row 27
$ws.Range("H27").Value = 18231.666
$ws.Range("I27").Value = 799.5
$ws.Range("J27").Value = 26947.75
$ws.Range("K27").Value = 799.5
$ws.Range("L27").Value = 26947.75
$ws.Range("M27").Value = -692.5
$ws.Range("N27").Value = -27161.75
# Row 46
$ws.Range("H46").Value = 1479.963
$ws.Range("I46").Value = 1176.4117
$ws.Range("J46").Value = 1996
$ws.Range("K46").Value = 1176.4117
$ws.Range("L46").Value = 1996
$ws.Range("M46").Value = -988.4117000000001
$ws.Range("N46").Value = -2372
# Row 61
$ws.Range("H61").Value = 9263.531000000001
$ws.Range("I61").Value = 11030.238
$ws.Range("J61").Value = 5890.727
$ws.Range("K61").Value = 11030.238
$ws.Range("L61").Value = 5890.727
$ws.Range("M61").Value = -10828.238
$ws.Range("N61").Value = -6294.727
# Row 75
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("N75").ClearContents()
# Row 78
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("N78").ClearContents()
# Row 113
$ws.Range("H113").Value = 9263.531000000001
$ws.Range("I113").Value = 11030.238
$ws.Range("J113").Value = 5890.727
$ws.Range("K113").Value = 11030.238
$ws.Range("L113").Value = 5890.727
$ws.Range("M113").Value = -8860.237999999999
$ws.Range("N113").Value = -10230.727
# Row 136
$ws.Range("H136").Value = 3607.0393
$ws.Range("I136").Value = 2161.7805
$ws.Range("K136").Value = 6485.3415
$ws.Range("M136").Value = -3935.3415

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 6173783
$ws.Range("I107").Value = 9260008
$ws.Range("K107").Value = 27780024
$ws.Range("M107").Value = -27778104
# Row 113
$ws.Range("H113").Value = 552.26666
$ws.Range("I113").Value = 400.57144
$ws.Range("J113").Value = 685
$ws.Range("K113").Value = 1201.71432
$ws.Range("L113").Value = 2055
$ws.Range("M113").Value = 968.28568
$ws.Range("N113").Value = -6395
